$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "correo" column (G) with header + three e-mail addresses.
$ws.Range("G1").Value = "correo"
$ws.Range("G2").Value = "mariavyeguezp@gmail.com"
$ws.Range("G3").Value = "r.gzlobos@gmail.com"
$ws.Range("G4").Value = "ma.yeguez@duocuc.cl"

# Turn the e-mail addresses into real (mailto:) hyperlinks - this is what
# drives Excel to auto-create the "Hyperlink" cell style / font as well.
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:mariavyeguezp@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:r.gzlobos@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:ma.yeguez@duocuc.cl")

# Match the column width Excel auto-applied to the new column.
$ws.Columns("G").ColumnWidth = 23.26953125

# Restore the view state captured in the saved workbook (scrolled so column
# C is left-most visible, with G10 selected/active).
$ws.Range("G10").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
